# Generate Report for Archive
#
# Updates the localization-status workbook so that the two files that have
# finished translation (0d74cade-2ec6-46bd-bee8-8f256f7b3fb6.md and
# 24ad6be8-8390-44f8-99d0-dff2e1b4ea7d.md) show a Status of "In Translation"
# instead of "Ready for handoff" on every sheet that tracks them.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: File Name in column A, per-locale status in B (zh-cn) and C (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B8").Value = $newStatus
$wsOverview.Range("C8").Value = $newStatus
$wsOverview.Range("B9").Value = $newStatus
$wsOverview.Range("C9").Value = $newStatus

# --- zh-cn sheet: Status is column C
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C8").Value = $newStatus
$wsZhCn.Range("C9").Value = $newStatus

# --- de-de sheet: Status is column C
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C8").Value = $newStatus
$wsDeDe.Range("C9").Value = $newStatus
